$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the "Purpose" row (A5, merged A5:E5): locial -> logical
$ws.Range("A5").Value = "Purpose:  Unit test the logical structure of the Export_Data Class and its Interface"

# Update the test-case row (row 9) for the new Input/Output stream null test
$ws.Range("B9").Value = "Input_Stream_And_Output_Stream_Null"
$ws.Range("C9").Value = "input_Stream = Null, output_Stream = Null"
$ws.Range("D9").Value = "NullPointerException"
$ws.Range("E9").Value = "NullPointerException thrown from the Input_Stream check in the resulting if statement check and cannot check if the Output_Stream exception throw works"

# Update the "Comments" row (A7, merged A7:E7) to mention FileOutputStream as well
$ws.Range("A7").Value = "Comments: Cannot fully test all pathways due to the requirement of a valid FileInputStream and FileOutputStream Object, both of which cannot be mocked due to Framework of Android Studio and it being a core object"

# B9 loses its distinct vertical-center / explicit-black-font styling and now shares the
# same left/vertical-center/wrap border style that C9:E9 already use
$ws.Range("C9").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 9 grows taller to fit the new, longer wrapped text
$ws.Rows.Item(9).RowHeight = 60

# Update the current selection to match the edited document (was A7:E7)
$ws.Range("A3:E6").Select()
